$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Row data: row number, date serial (yyyy-mm-dd), Low (B), High (C)
# Dates converted from inline-string "yyyy-mm-dd" text to real Excel date
# serials with a yyyy-mm-dd number format. Rows 87-102 (the older
# pre-split Apple prices) are re-adjusted for the Aug 2020 4-for-1 stock
# split (divide by 4).
$data = @(
    @(3, 44193, 137.34, 133.51),
    @(4, 44189, 133.46, 131.1),
    @(5, 44188, 132.43, 130.78),
    @(6, 44187, 134.405, 129.65),
    @(7, 44186, 128.31, 123.449),
    @(8, 44183, 129.1, 126.12),
    @(9, 44182, 129.58, 128.045),
    @(10, 44181, 128.37, 126.56),
    @(11, 44180, 127.9, 124.13),
    @(12, 44179, 123.35, 121.54),
    @(13, 44176, 122.76, 120.55),
    @(14, 44175, 123.87, 120.15),
    @(15, 44174, 125.95, 121),
    @(16, 44173, 124.98, 123.09),
    @(17, 44172, 124.57, 122.25),
    @(18, 44169, 122.8608, 121.52),
    @(19, 44168, 123.78, 122.21),
    @(20, 44167, 123.37, 120.89),
    @(21, 44166, 123.4693, 120.01),
    @(22, 44165, 120.97, 116.81),
    @(23, 44162, 117.49, 116.22),
    @(24, 44160, 116.75, 115.17),
    @(25, 44159, 115.85, 112.59),
    @(26, 44158, 117.6202, 113.75),
    @(27, 44155, 118.77, 117.29),
    @(28, 44154, 119.06, 116.81),
    @(29, 44153, 119.82, 118),
    @(30, 44152, 120.6741, 118.96),
    @(31, 44151, 120.99, 118.146),
    @(32, 44148, 119.6717, 117.87),
    @(33, 44147, 120.53, 118.57),
    @(34, 44146, 119.63, 116.44),
    @(35, 44145, 117.59, 114.13),
    @(36, 44144, 121.99, 116.05),
    @(37, 44141, 119.2, 116.13),
    @(38, 44140, 119.62, 116.8686),
    @(39, 44139, 115.59, 112.35),
    @(40, 44138, 111.49, 108.73),
    @(41, 44137, 110.68, 107.32),
    @(42, 44134, 111.99, 107.72),
    @(43, 44133, 116.93, 112.2),
    @(44, 44132, 115.43, 111.1),
    @(45, 44131, 117.28, 114.5399),
    @(46, 44130, 116.55, 112.88),
    @(47, 44127, 116.55, 114.28),
    @(48, 44126, 118.04, 114.59),
    @(49, 44125, 118.705, 116.45),
    @(50, 44124, 118.98, 115.63),
    @(51, 44123, 120.419, 115.66),
    @(52, 44120, 121.548, 118.81),
    @(53, 44119, 121.2, 118.15),
    @(54, 44118, 123.03, 119.62),
    @(55, 44117, 125.39, 119.65),
    @(56, 44116, 125.18, 119.2845),
    @(57, 44113, 117, 114.92),
    @(58, 44112, 116.4, 114.5901),
    @(59, 44111, 115.55, 114.13),
    @(60, 44110, 116.12, 112.25),
    @(61, 44109, 116.65, 113.55),
    @(62, 44106, 115.37, 112.22),
    @(63, 44105, 117.72, 115.83),
    @(64, 44104, 117.26, 113.62),
    @(65, 44103, 115.31, 113.57),
    @(66, 44102, 115.32, 112.78),
    @(67, 44099, 112.44, 107.67),
    @(68, 44098, 110.25, 105),
    @(69, 44097, 112.11, 106.77),
    @(70, 44096, 112.86, 109.16),
    @(71, 44095, 110.19, 103.1),
    @(72, 44092, 110.88, 106.09),
    @(73, 44091, 112.2, 108.71),
    @(74, 44090, 116, 112.04),
    @(75, 44089, 118.829, 113.61),
    @(76, 44088, 115.93, 112.8),
    @(77, 44085, 115.23, 110),
    @(78, 44084, 120.5, 112.5),
    @(79, 44083, 119.14, 115.26),
    @(80, 44082, 118.99, 112.68),
    @(81, 44078, 123.7, 110.89),
    @(82, 44077, 128.84, 120.5),
    @(83, 44076, 137.98, 127),
    @(84, 44075, 134.8, 130.53),
    @(85, 44074, 131, 126),
    @(86, 44071, 505.77, 498.31),
    @(87, 44070, 127.485, 123.8325),
    @(88, 44069, 126.9925, 125.0825),
    @(89, 44068, 125.1793, 123.0525),
    @(90, 44067, 128.785, 123.93625),
    @(91, 44064, 124.868, 119.25),
    @(92, 44063, 118.392, 115.733375),
    @(93, 44062, 117.1625, 115.61),
    @(94, 44061, 116, 114.0075),
    @(95, 44060, 116.0875, 113.962525),
    @(96, 44057, 115, 113.045),
    @(97, 44056, 116.0425, 113.9275),
    @(98, 44055, 113.275, 110.2975),
    @(99, 44054, 112.4825, 109.106675),
    @(100, 44053, 113.775, 110),
    @(101, 44050, 113.675, 110.2925),
    @(102, 44049, 114.4125, 109.7975)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Fix the chart's series references: they were off by one column
# ('Summary'!B/C/D -> 'Summary'!A/B/C) and needed a wider range to cover
# the full data set (rows 2-100 instead of 2-7).
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$s1 = $chart.SeriesCollection().Item(1)
$s1.Formula = "=SERIES('Summary'!`$A`$1,,'Summary'!`$A`$2:`$A`$100,1)"

$s2 = $chart.SeriesCollection().Item(2)
$s2.Formula = "=SERIES('Summary'!`$B`$1,,'Summary'!`$B`$2:`$B`$100,2)"

$s3 = $chart.SeriesCollection().Item(3)
$s3.Formula = "=SERIES('Summary'!`$C`$1,,'Summary'!`$C`$2:`$C`$100,3)"
